$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = 44245
$ws.Range("B20").Value = 8
$ws.Range("D20").Value = "Meeting+Development WebUI"

$ws.Range("F19").Select()
